$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "0.100" or
# "58.638.43" round-trip exactly instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.638.43"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.587.87"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("D5").Value = "552.38"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").Value = "143.06"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +4.68%  "
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("E11").Value = "  +4.27%  "
$ws.Range("D12").Value = "0.334"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "3.047.21"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "58.569.68"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "20.77"
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "2.596.99"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").Value = "4.43"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "335.89"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "10.02"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "6.11"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "66.24"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").Value = "0.426"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -4.06%  "
$ws.Range("D27").Value = "7.10"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").Value = "0.0₃0750"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "5.90"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "153.55"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").Value = "18.86"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "3.90"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").Value = "0.866"
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("D36").Value = "37.11"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "0.818"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").Value = "3.59"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "281.85"
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "10.61"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "0.0529"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").Value = "0.0226"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "1.908.22"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "4.41"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "114.65"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "17.69"
$ws.Range("E51").Value = "  -3.23%  "
